# Scheduled market-price refresh: update currentAveragePrice/NQ/HQ and
# derived Leve profit columns (H,I,J,K,L,M,N) across the profession sheets.
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 11
$ws_ALC.Range("H11").Value = 47.642857
$ws_ALC.Range("I11").Value = 47.642857
$ws_ALC.Range("K11").Value = 47.642857
$ws_ALC.Range("M11").Value = 92.35714300000001

# ALC row 28
$ws_ALC.Range("H28").Value = 812.875
$ws_ALC.Range("I28").Value = 538.7
$ws_ALC.Range("K28").Value = 538.7
$ws_ALC.Range("M28").Value = -53.70000000000005

# ALC row 54
$ws_ALC.Range("H54").Value = 5000
$ws_ALC.Range("I54").Value = 5000
$ws_ALC.Range("K54").Value = 5000
$ws_ALC.Range("M54").Value = -4514

# ALC row 86
$ws_ALC.Range("H86").Value = 2751.5
$ws_ALC.Range("I86").Value = 2092.1333
$ws_ALC.Range("J86").Value = 3650.6365
$ws_ALC.Range("K86").Value = 2092.1333
$ws_ALC.Range("L86").Value = 3650.6365
$ws_ALC.Range("M86").Value = -969.1333
$ws_ALC.Range("N86").Value = -5896.636500000001

# ALC row 89
$ws_ALC.Range("H89").Value = 2751.5
$ws_ALC.Range("I89").Value = 2092.1333
$ws_ALC.Range("J89").Value = 3650.6365
$ws_ALC.Range("K89").Value = 10460.6665
$ws_ALC.Range("L89").Value = 18253.1825
$ws_ALC.Range("M89").Value = -4844.666499999999
$ws_ALC.Range("N89").Value = -29485.1825

# ALC row 132
$ws_ALC.Range("H132").Value = 10783.049
$ws_ALC.Range("I132").Value = 1656.3773
$ws_ALC.Range("K132").Value = 4969.1319
$ws_ALC.Range("M132").Value = -2439.1319

# ARM row 32
$ws_ARM.Range("H32").Value = 15157624
$ws_ARM.Range("J32").Value = 4666.3335
$ws_ARM.Range("L32").Value = 4666.3335
$ws_ARM.Range("N32").Value = -5240.3335

# ARM row 61
$ws_ARM.Range("H61").Value = 2641.45
$ws_ARM.Range("I61").Value = 2712.7222
$ws_ARM.Range("K61").Value = 2712.7222
$ws_ARM.Range("M61").Value = -2500.7222

# ARM row 132
$ws_ARM.Range("H132").Value = 3109
$ws_ARM.Range("I132").Value = 3109
$ws_ARM.Range("K132").Value = 9327
$ws_ARM.Range("M132").Value = -6797

# ARM row 136
$ws_ARM.Range("H136").Value = 2641.45
$ws_ARM.Range("I136").Value = 2712.7222
$ws_ARM.Range("K136").Value = 8138.1666
$ws_ARM.Range("M136").Value = -5588.1666

# BSM row 36
$ws_BSM.Range("H36").Value = 1666
$ws_BSM.Range("I36").Value = 1666
$ws_BSM.Range("K36").Value = 1666
$ws_BSM.Range("M36").Value = -1132

# BSM row 105
$ws_BSM.Range("H105").Value = 2712.1333
$ws_BSM.Range("I105").Value = 1517
$ws_BSM.Range("K105").Value = 1517
$ws_BSM.Range("M105").Value = 230

# CRP row 58
$ws_CRP.Range("H58").Value = 1830.6111
$ws_CRP.Range("I58").Value = 952.92
$ws_CRP.Range("J58").Value = 3825.3635
$ws_CRP.Range("K58").Value = 952.92
$ws_CRP.Range("L58").Value = 3825.3635
$ws_CRP.Range("M58").Value = -749.92
$ws_CRP.Range("N58").Value = -4231.363499999999

# CRP row 99
$ws_CRP.Range("H99").Value = 22444636
$ws_CRP.Range("I99").Value = 4067222.8
$ws_CRP.Range("K99").Value = 4067222.8
$ws_CRP.Range("M99").Value = -4065724.8

# CRP row 118
$ws_CRP.Range("H118").Value = 49545.453
$ws_CRP.Range("J118").Value = 49545.453
$ws_CRP.Range("L118").Value = 49545.453
$ws_CRP.Range("N118").Value = -52859.453

# CRP row 122
$ws_CRP.Range("H122").Value = 467669.28
$ws_CRP.Range("I122").Value = 730337.4399999999
$ws_CRP.Range("K122").Value = 2191012.32
$ws_CRP.Range("M122").Value = -2188562.32

# CRP row 126
$ws_CRP.Range("H126").Value = 22444636
$ws_CRP.Range("I126").Value = 4067222.8
$ws_CRP.Range("K126").Value = 12201668.4
$ws_CRP.Range("M126").Value = -12199198.4

# CRP row 134
$ws_CRP.Range("H134").Value = 2020.6364
$ws_CRP.Range("I134").Value = 2165.8928
$ws_CRP.Range("J134").Value = 1207.2
$ws_CRP.Range("K134").Value = 6497.678400000001
$ws_CRP.Range("L134").Value = 3621.6
$ws_CRP.Range("M134").Value = -3962.678400000001
$ws_CRP.Range("N134").Value = -8691.6

# CRP row 136
$ws_CRP.Range("H136").Value = 1830.6111
$ws_CRP.Range("I136").Value = 952.92
$ws_CRP.Range("J136").Value = 3825.3635
$ws_CRP.Range("K136").Value = 2858.76
$ws_CRP.Range("L136").Value = 11476.0905
$ws_CRP.Range("M136").Value = -308.7599999999998
$ws_CRP.Range("N136").Value = -16576.0905

# CUL row 34
$ws_CUL.Range("H34").Value = 98
$ws_CUL.Range("J34").Value = 0
$ws_CUL.Range("L34").Value = 0
$ws_CUL.Range("N34").ClearContents()  # cell removed entirely in target

# CUL row 61
$ws_CUL.Range("H61").Value = 208.41176
$ws_CUL.Range("I61").Value = 178.61539
$ws_CUL.Range("J61").Value = 305.25
$ws_CUL.Range("K61").Value = 535.84617
$ws_CUL.Range("L61").Value = 915.75
$ws_CUL.Range("M61").Value = -320.84617
$ws_CUL.Range("N61").Value = -1345.75

# CUL row 129
$ws_CUL.Range("H129").Value = 1739.6
$ws_CUL.Range("I129").Value = 993.625
$ws_CUL.Range("K129").Value = 2980.875
$ws_CUL.Range("M129").Value = 2019.125

# CUL row 131
$ws_CUL.Range("H131").Value = 3570.9473
$ws_CUL.Range("J131").Value = 4970.222
$ws_CUL.Range("L131").Value = 14910.666
$ws_CUL.Range("N131").Value = -24990.666

# CUL row 137
$ws_CUL.Range("H137").Value = 2762.56
$ws_CUL.Range("I137").Value = 2370.3076
$ws_CUL.Range("J137").Value = 3187.5
$ws_CUL.Range("K137").Value = 7110.9228
$ws_CUL.Range("L137").Value = 9562.5
$ws_CUL.Range("M137").Value = -2010.9228
$ws_CUL.Range("N137").Value = -19762.5

# GSM row 80
$ws_GSM.Range("H80").Value = 69918.89999999999
$ws_GSM.Range("I80").Value = 111283.82
$ws_GSM.Range("J80").Value = 13042.125
$ws_GSM.Range("K80").Value = 111283.82
$ws_GSM.Range("L80").Value = 13042.125
$ws_GSM.Range("M80").Value = -110285.82
$ws_GSM.Range("N80").Value = -15038.125

# GSM row 83
$ws_GSM.Range("H83").Value = 69918.89999999999
$ws_GSM.Range("I83").Value = 111283.82
$ws_GSM.Range("J83").Value = 13042.125
$ws_GSM.Range("K83").Value = 556419.1000000001
$ws_GSM.Range("L83").Value = 65210.625
$ws_GSM.Range("M83").Value = -551427.1000000001
$ws_GSM.Range("N83").Value = -75194.625

# GSM row 102
$ws_GSM.Range("H102").Value = 10999.5
$ws_GSM.Range("J102").Value = 13332.667
$ws_GSM.Range("L102").Value = 13332.667
$ws_GSM.Range("N102").Value = -16576.667

# GSM row 113
$ws_GSM.Range("H113").Value = 1755.8334
$ws_GSM.Range("I113").Value = 1179
$ws_GSM.Range("J113").Value = 3774.75
$ws_GSM.Range("K113").Value = 1179
$ws_GSM.Range("L113").Value = 3774.75
$ws_GSM.Range("M113").Value = 991
$ws_GSM.Range("N113").Value = -8114.75

# GSM row 122
$ws_GSM.Range("H122").Value = 5610.857
$ws_GSM.Range("I122").Value = 6233.6
$ws_GSM.Range("K122").Value = 18700.8
$ws_GSM.Range("M122").Value = -16250.8

# LTW row 132
$ws_LTW.Range("H132").Value = 5359.1816
$ws_LTW.Range("I132").Value = 3349.5
$ws_LTW.Range("K132").Value = 10048.5
$ws_LTW.Range("M132").Value = -7518.5

# WVR row 81
$ws_WVR.Range("H81").Value = 12350489
$ws_WVR.Range("I81").Value = 4573.75
$ws_WVR.Range("J81").Value = 22227222
$ws_WVR.Range("K81").Value = 9147.5
$ws_WVR.Range("L81").Value = 44454444
$ws_WVR.Range("M81").Value = -8086.5
$ws_WVR.Range("N81").Value = -44456566

# WVR row 84
$ws_WVR.Range("H84").Value = 12350489
$ws_WVR.Range("I84").Value = 4573.75
$ws_WVR.Range("J84").Value = 22227222
$ws_WVR.Range("K84").Value = 45737.5
$ws_WVR.Range("L84").Value = 222272220
$ws_WVR.Range("M84").Value = -40433.5
$ws_WVR.Range("N84").Value = -222282828

# WVR row 136
$ws_WVR.Range("H136").Value = 3449.5
$ws_WVR.Range("I136").Value = 1478.8
$ws_WVR.Range("K136").Value = 4436.4
$ws_WVR.Range("M136").Value = -1886.4

# WVR row 141
$ws_WVR.Range("H141").Value = 86247.88
$ws_WVR.Range("J141").Value = 86013.375
$ws_WVR.Range("L141").Value = 86013.375
$ws_WVR.Range("N141").Value = -96373.375
